$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 72

# Write the date as a formula-derived text string and convert it in place
# to a plain value; this avoids Excel auto-detecting "2020-08-10" as a
# date literal (which would store it as a numeric serial + date format)
# and keeps the cell free of any new cell style.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.Formula = "=""2020-08-10"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

$ws.Cells.Item($row, 2).Value = 485836
$ws.Cells.Item($row, 3).Value = 532028
$ws.Cells.Item($row, 4).Value = 79213
$ws.Cells.Item($row, 5).Value = 53003
$ws.Cells.Item($row, 6).Value = 26.56
